$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the existing _GoBack bookmark (it currently sits right
# after "No additional discoveries..." at the end of that paragraph).
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------------
# Step 2: insert the two new whole paragraphs right after "UAT Run"
# and before the (currently empty) trailing paragraph.
# ---------------------------------------------------------------------
$pTrailing = $d.Paragraphs.Item(136)
$rTrailing = $pTrailing.Range
$rTrailing.Collapse(1)
$rTrailing.InsertBefore("The run has passed, but I have noted that the UI at present doesn" + [char]0x2019 + "t allow for a single run, and even if it did, it would be difficult to test a random game of chance in a single run-through of the UAT as recorded.  " + "`r")

$pTrailing2 = $d.Paragraphs.Item(137)
$rTrailing2 = $pTrailing2.Range
$rTrailing2.Collapse(1)
$rTrailing2.InsertBefore("I have chosen an appropriate part of the log to highlight the successful resolution of the bug." + "`r")

# ---------------------------------------------------------------------
# Step 3: fill in the final paragraph (still empty) with the full
# sentence as a single run first (the trailing CR trick makes the new
# run inherit the paragraph mark's rPr, i.e. lang=en-US); we will then
# split it into the required runs.
# ---------------------------------------------------------------------
$pFinal = $d.Paragraphs.Item(138)
$rFinal = $pFinal.Range
$rFinal.Collapse(1)
$finalStart = $rFinal.Start
$run1 = "I have also recorded t"
$run2 = "he " + [char]0x201C + "Breaking the bank" + [char]0x201D + " issue"
$run3 = " on the UAT."
$rFinal.InsertBefore($run1 + $run2 + $run3 + "`r")

# ---------------------------------------------------------------------
# Step 4: split off run1 | (run2+run3) using a throw-away bookmark,
# then delete that bookmark -- the run split survives the deletion.
# ---------------------------------------------------------------------
$splitPos1 = $finalStart + $run1.Length
$rSplit1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("TempSplit1", $rSplit1)
$tempBm1 = $d.Bookmarks.Item("TempSplit1")
$tempBm1.Delete()

# ---------------------------------------------------------------------
# Step 5: split off run2 | run3 with the real _GoBack bookmark (this
# one stays, so it naturally creates the two surrounding runs).
# ---------------------------------------------------------------------
$splitPos2 = $finalStart + $run1.Length + $run2.Length
$rSplit2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $rSplit2)

Write-Host "Done"
